# Gyroscope data update (changes from May 9th):
#  - 9 new sample rows are inserted right after the header row, pushing the
#    20 existing data rows down by 9 rows.
#  - 1 new sample row is appended after the (now shifted) last row.
# Net effect: sheet grows from A1:C21 to A1:C31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Preserve the current data block (rows 2-21, columns A-C) before we
#    overwrite anything, then write it back 9 rows further down (rows 11-30).
#    Value2 round-trips numbers cleanly (no locale-formatted strings).
$existingData = $ws.Range("A2:C21").Value2
$ws.Range("A11:C30").Value2 = $existingData

# 2) Fill the newly freed rows 2-10 with the 9 new samples.
$newTopRows = @(
    @(-0.0114537235349416, -0.0096211275085806, -0.0482583530247211),
    @(0.0445931628346443, 0.1122464910149574, -0.0378736443817615),
    @(0.0612392425537109, 0.09758572280406951, -0.0021380283869802),
    @(0.0088575463742017, 0.1237002089619636, 0.0548251569271087),
    @(-0.0221438650041818, 0.0061086523346602, 0.0325285755097866),
    @(0.0332921557128429, -0.0615446716547012, 0.093156948685646),
    @(-0.4489859640598297, -1.353219270706177, 0.4497495293617248),
    @(-0.3888157308101654, -3.63083028793335, -0.1369865238666534),
    @(-0.6565274000167847, -2.371837139129639, 0.1600466966629028)
)

for ($i = 0; $i -lt $newTopRows.Count; $i++) {
    $r = 2 + $i
    $row = $newTopRows[$i]
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
}

# 3) Append the one new sample row at the end (row 31).
$ws.Cells.Item(31, 1).Value2 = 0.0546724386513233
$ws.Cells.Item(31, 2).Value2 = -0.1007927656173706
$ws.Cells.Item(31, 3).Value2 = 0.2141082733869552

Write-Host "done"
